$wb = $excel.ActiveWorkbook
$wsCombined = $wb.Worksheets.Item("Combined")
$wsKuCoin = $wb.Worksheets.Item("KuCoin")
$wsBinance = $wb.Worksheets.Item("Binance")

# --- Combined sheet: refreshed spread/slippage figures + re-ordered trading pairs ---
$wsCombined.Range("B2").Value = "BTC-USDT"
$wsCombined.Range("C2").Value = 0.0002355707053659351
$wsCombined.Range("D2").Value = 0.01154296456310222
$wsCombined.Range("E2").Value = 0.01177856301531213
$wsCombined.Range("B3").Value = "ADA-USDT"
$wsCombined.Range("C3").Value = 0.01292478406153876
$wsCombined.Range("D3").Value = 0.01339764201500563
$wsCombined.Range("E3").Value = 0.02632582863304762
$wsCombined.Range("B4").Value = "ETH-BTC"
$wsCombined.Range("C4").Value = 0.001883451990810638
$wsCombined.Range("D4").Value = -0.001883451990810638
$wsCombined.Range("E4").Value = 0
$wsCombined.Range("B5").Value = "LINK-USDT"
$wsCombined.Range("C5").Value = 0.003358950932448884
$wsCombined.Range("D5").Value = 0.05777395603804681
$wsCombined.Range("E5").Value = 0.06113496046382026
$wsCombined.Range("B6").Value = "AVAX-USDT"
$wsCombined.Range("C6").Value = 0.01247411620885526
$wsCombined.Range("D6").Value = -0.02245340917596074
$wsCombined.Range("E6").Value = -0.009980537951007755
$wsCombined.Range("B7").Value = "RUNE-USDT"
$wsCombined.Range("C7").Value = 0.0269079412061455
$wsCombined.Range("D7").Value = -0.04372540445999391
$wsCombined.Range("E7").Value = -0.01682198970494792
$wsCombined.Range("B8").Value = "BONK-USDT"
$wsCombined.Range("C8").Value = 0.1621304632791797
$wsCombined.Range("D8").Value = -0.0137983373003505
$wsCombined.Range("E8").Value = 0.1485730080851455
$wsCombined.Range("B9").Value = "ETH-USDT"
$wsCombined.Range("C9").Value = 0.0004437344692830818
$wsCombined.Range("D9").Value = -0.0004437344692830818
$wsCombined.Range("E9").Value = 0
$wsCombined.Range("B10").Value = "DOT-USDT"
$wsCombined.Range("C10").Value = 0.01380357512596223
$wsCombined.Range("D10").Value = 0.09800538339429139
$wsCombined.Range("E10").Value = 0.1118243942845278
$wsCombined.Range("B11").Value = "SOL-USDT"
$wsCombined.Range("C11").Value = 0.003926958570587229
$wsCombined.Range("D11").Value = 0
$wsCombined.Range("E11").Value = 0.003927112786679382
$wsCombined.Range("B12").Value = "DOGE-USDT"
$wsCombined.Range("C12").Value = 0.02101281781887595
$wsCombined.Range("D12").Value = 0
$wsCombined.Range("E12").Value = 0.02101723413199467
$wsCombined.Range("B13").Value = "ATOM-USDT"
$wsCombined.Range("C13").Value = 0.004472992073863511
$wsCombined.Range("D13").Value = 0.01431357463634417
$wsCombined.Range("E13").Value = 0.01878740706943678
$wsCombined.Range("B14").Value = "ETH-USDC"
$wsCombined.Range("C14").Value = 0.002661910719512047
$wsCombined.Range("D14").Value = 0.04658343759151125
$wsCombined.Range("E14").Value = 0.04924665921312389
$wsCombined.Range("B15").Value = "MATIC-USDT"
$wsCombined.Range("C15").Value = 0.01161575095829817
$wsCombined.Range("D15").Value = -0.01161575095829817
$wsCombined.Range("B16").Value = "INJ-USDT"
$wsCombined.Range("C16").Value = 0.08347781416554229
$wsCombined.Range("D16").Value = -0.1091632954472564
$wsCombined.Range("E16").Value = -0.02570694087404457
$wsCombined.Range("B17").Value = "TIA-USDT"
$wsCombined.Range("C17").Value = 0.08548769969392764
$wsCombined.Range("D17").Value = 0
$wsCombined.Range("E17").Value = 0.08556084369103784
$wsCombined.Range("B18").Value = "JTO-USDT"
$wsCombined.Range("C18").Value = 0.1519401589527811
$wsCombined.Range("D18").Value = 0.04517139860758508
$wsCombined.Range("E18").Value = 0.1974115049147247
$wsCombined.Range("B19").Value = "XRP-USDT"
$wsCombined.Range("C19").Value = 0.001601819667152367
$wsCombined.Range("D19").Value = 0.05766550801711175
$wsCombined.Range("E19").Value = 0.05926827705518237
$wsCombined.Range("B20").Value = "FET-USDT"
$wsCombined.Range("C20").Value = 0.01381215469613108
$wsCombined.Range("D20").Value = 0.1381215469613261
$wsCombined.Range("E20").Value = 0.1519546898742906
$wsCombined.Range("B21").Value = "BTC-USDC"
$wsCombined.Range("C21").Value = 0.01012560459277558
$wsCombined.Range("D21").Value = 0.00777081282702856
$wsCombined.Range("E21").Value = 0.01789822972377508
$wsCombined.Range("C22").Value = 0.00002355197178765182
$wsCombined.Range("D22").Value = -0.00002355197178765182
$wsCombined.Range("B23").Value = "SOLUSDT"
$wsCombined.Range("C23").Value = 0.01307702366941953
$wsCombined.Range("D23").Value = -0.01307702366941953
$wsCombined.Range("E23").Value = 0
$wsCombined.Range("B24").Value = "FETUSDT"
$wsCombined.Range("C24").Value = 0.02762049440684684
$wsCombined.Range("D24").Value = -0.04143074161027026
$wsCombined.Range("E24").Value = -0.01381406271584321
$wsCombined.Range("C25").Value = 0.01905938140871321
$wsCombined.Range("D25").Value = -0.03659401230473261
$wsCombined.Range("E25").Value = -0.01753797352528494
$wsCombined.Range("B26").Value = "BONKUSDT"
$wsCombined.Range("C26").Value = 0.03441156228492158
$wsCombined.Range("D26").Value = 0
$wsCombined.Range("E26").Value = 0.03442340791737766
$wsCombined.Range("B27").Value = "ADAUSDT"
$wsCombined.Range("C27").Value = 0.01575795776868894
$wsCombined.Range("D27").Value = 0
$wsCombined.Range("E27").Value = 0.01576044129237195
$wsCombined.Range("B28").Value = "JTOUSDT"
$wsCombined.Range("C28").Value = 0.02210921954455465
$wsCombined.Range("E28").Value = 0.02211410880141988
$wsCombined.Range("B29").Value = "ETHUSDC"
$wsCombined.Range("C29").Value = 0.001330595263976798
$wsCombined.Range("D29").Value = -0.001330595263976798
$wsCombined.Range("B30").Value = "ETHUSDT"
$wsCombined.Range("C30").Value = 0.0004435750373607078
$wsCombined.Range("D30").Value = -0.0004435750373607078
$wsCombined.Range("E30").Value = 0
$wsCombined.Range("B31").Value = "RUNEUSDT"
$wsCombined.Range("C31").Value = 0.01681520094165687
$wsCombined.Range("E31").Value = 0.01681802892701537
$wsCombined.Range("B32").Value = "MATICUSDT"
$wsCombined.Range("C32").Value = 0.01161710037174593
$wsCombined.Range("D32").Value = 0
$wsCombined.Range("E32").Value = 0.01161845009875555
$wsCombined.Range("B33").Value = "DOTUSDT"
$wsCombined.Range("C33").Value = 0.01379310344828047
$wsCombined.Range("D33").Value = -0.04137931034482915
$wsCombined.Range("E33").Value = -0.02759001241550255
$wsCombined.Range("B34").Value = "ATOMUSDT"
$wsCombined.Range("C34").Value = 0.008942144326204469
$wsCombined.Range("D34").Value = -0.02682643297862929
$wsCombined.Range("E34").Value = -0.01788588803434688
$wsCombined.Range("B35").Value = "ETHBTC"
$wsCombined.Range("C35").Value = 0.01883239171375341
$wsCombined.Range("E35").Value = 0.0188359389715635
$wsCombined.Range("B36").Value = "LINKUSDT"
$wsCombined.Range("C36").Value = 0.00671546571754379
$wsCombined.Range("D36").Value = -0.00671546571754379
$wsCombined.Range("B37").Value = "DOGEUSDT"
$wsCombined.Range("C37").Value = 0.01049979000421042
$wsCombined.Range("D37").Value = -0.01049979000421042
$wsCombined.Range("B38").Value = "XRPUSDT"
$wsCombined.Range("C38").Value = 0.01601537475976762
$wsCombined.Range("D38").Value = -0.01601537475976762
$wsCombined.Range("B39").Value = "BTCUSDC"
$wsCombined.Range("C39").Value = 0.00004709817543269348
$wsCombined.Range("D39").Value = -0.007370864456390223
$wsCombined.Range("E39").Value = -0.007323769730319446
$wsCombined.Range("B40").Value = "AVAXUSDT"
$wsCombined.Range("C40").Value = 0.02495632642874472
$wsCombined.Range("D40").Value = 0
$wsCombined.Range("E40").Value = 0.0249625561657464
$wsCombined.Range("B41").Value = "INJUSDT"
$wsCombined.Range("C41").Value = 0.009631127805066338
$wsCombined.Range("D41").Value = -0.006420751870051829
$wsCombined.Range("E41").Value = 0.003210685160205707

# --- KuCoin sheet: mirrors Combined rows 2-21 ---
$wsKuCoin.Range("B2").Value = "BTC-USDT"
$wsKuCoin.Range("C2").Value = 0.0002355707053659351
$wsKuCoin.Range("D2").Value = 0.01154296456310222
$wsKuCoin.Range("E2").Value = 0.01177856301531213
$wsKuCoin.Range("B3").Value = "ADA-USDT"
$wsKuCoin.Range("C3").Value = 0.01292478406153876
$wsKuCoin.Range("D3").Value = 0.01339764201500563
$wsKuCoin.Range("E3").Value = 0.02632582863304762
$wsKuCoin.Range("B4").Value = "ETH-BTC"
$wsKuCoin.Range("C4").Value = 0.001883451990810638
$wsKuCoin.Range("D4").Value = -0.001883451990810638
$wsKuCoin.Range("E4").Value = 0
$wsKuCoin.Range("B5").Value = "LINK-USDT"
$wsKuCoin.Range("C5").Value = 0.003358950932448884
$wsKuCoin.Range("D5").Value = 0.05777395603804681
$wsKuCoin.Range("E5").Value = 0.06113496046382026
$wsKuCoin.Range("B6").Value = "AVAX-USDT"
$wsKuCoin.Range("C6").Value = 0.01247411620885526
$wsKuCoin.Range("D6").Value = -0.02245340917596074
$wsKuCoin.Range("E6").Value = -0.009980537951007755
$wsKuCoin.Range("B7").Value = "RUNE-USDT"
$wsKuCoin.Range("C7").Value = 0.0269079412061455
$wsKuCoin.Range("D7").Value = -0.04372540445999391
$wsKuCoin.Range("E7").Value = -0.01682198970494792
$wsKuCoin.Range("B8").Value = "BONK-USDT"
$wsKuCoin.Range("C8").Value = 0.1621304632791797
$wsKuCoin.Range("D8").Value = -0.0137983373003505
$wsKuCoin.Range("E8").Value = 0.1485730080851455
$wsKuCoin.Range("B9").Value = "ETH-USDT"
$wsKuCoin.Range("C9").Value = 0.0004437344692830818
$wsKuCoin.Range("D9").Value = -0.0004437344692830818
$wsKuCoin.Range("E9").Value = 0
$wsKuCoin.Range("B10").Value = "DOT-USDT"
$wsKuCoin.Range("C10").Value = 0.01380357512596223
$wsKuCoin.Range("D10").Value = 0.09800538339429139
$wsKuCoin.Range("E10").Value = 0.1118243942845278
$wsKuCoin.Range("B11").Value = "SOL-USDT"
$wsKuCoin.Range("C11").Value = 0.003926958570587229
$wsKuCoin.Range("D11").Value = 0
$wsKuCoin.Range("E11").Value = 0.003927112786679382
$wsKuCoin.Range("B12").Value = "DOGE-USDT"
$wsKuCoin.Range("C12").Value = 0.02101281781887595
$wsKuCoin.Range("D12").Value = 0
$wsKuCoin.Range("E12").Value = 0.02101723413199467
$wsKuCoin.Range("B13").Value = "ATOM-USDT"
$wsKuCoin.Range("C13").Value = 0.004472992073863511
$wsKuCoin.Range("D13").Value = 0.01431357463634417
$wsKuCoin.Range("E13").Value = 0.01878740706943678
$wsKuCoin.Range("B14").Value = "ETH-USDC"
$wsKuCoin.Range("C14").Value = 0.002661910719512047
$wsKuCoin.Range("D14").Value = 0.04658343759151125
$wsKuCoin.Range("E14").Value = 0.04924665921312389
$wsKuCoin.Range("B15").Value = "MATIC-USDT"
$wsKuCoin.Range("C15").Value = 0.01161575095829817
$wsKuCoin.Range("D15").Value = -0.01161575095829817
$wsKuCoin.Range("B16").Value = "INJ-USDT"
$wsKuCoin.Range("C16").Value = 0.08347781416554229
$wsKuCoin.Range("D16").Value = -0.1091632954472564
$wsKuCoin.Range("E16").Value = -0.02570694087404457
$wsKuCoin.Range("B17").Value = "TIA-USDT"
$wsKuCoin.Range("C17").Value = 0.08548769969392764
$wsKuCoin.Range("D17").Value = 0
$wsKuCoin.Range("E17").Value = 0.08556084369103784
$wsKuCoin.Range("B18").Value = "JTO-USDT"
$wsKuCoin.Range("C18").Value = 0.1519401589527811
$wsKuCoin.Range("D18").Value = 0.04517139860758508
$wsKuCoin.Range("E18").Value = 0.1974115049147247
$wsKuCoin.Range("B19").Value = "XRP-USDT"
$wsKuCoin.Range("C19").Value = 0.001601819667152367
$wsKuCoin.Range("D19").Value = 0.05766550801711175
$wsKuCoin.Range("E19").Value = 0.05926827705518237
$wsKuCoin.Range("B20").Value = "FET-USDT"
$wsKuCoin.Range("C20").Value = 0.01381215469613108
$wsKuCoin.Range("D20").Value = 0.1381215469613261
$wsKuCoin.Range("E20").Value = 0.1519546898742906
$wsKuCoin.Range("B21").Value = "BTC-USDC"
$wsKuCoin.Range("C21").Value = 0.01012560459277558
$wsKuCoin.Range("D21").Value = 0.00777081282702856
$wsKuCoin.Range("E21").Value = 0.01789822972377508

# --- Binance sheet: mirrors Combined rows 22-41 ---
$wsBinance.Range("C2").Value = 0.00002355197178765182
$wsBinance.Range("D2").Value = -0.00002355197178765182
$wsBinance.Range("B3").Value = "SOLUSDT"
$wsBinance.Range("C3").Value = 0.01307702366941953
$wsBinance.Range("D3").Value = -0.01307702366941953
$wsBinance.Range("E3").Value = 0
$wsBinance.Range("B4").Value = "FETUSDT"
$wsBinance.Range("C4").Value = 0.02762049440684684
$wsBinance.Range("D4").Value = -0.04143074161027026
$wsBinance.Range("E4").Value = -0.01381406271584321
$wsBinance.Range("C5").Value = 0.01905938140871321
$wsBinance.Range("D5").Value = -0.03659401230473261
$wsBinance.Range("E5").Value = -0.01753797352528494
$wsBinance.Range("B6").Value = "BONKUSDT"
$wsBinance.Range("C6").Value = 0.03441156228492158
$wsBinance.Range("D6").Value = 0
$wsBinance.Range("E6").Value = 0.03442340791737766
$wsBinance.Range("B7").Value = "ADAUSDT"
$wsBinance.Range("C7").Value = 0.01575795776868894
$wsBinance.Range("D7").Value = 0
$wsBinance.Range("E7").Value = 0.01576044129237195
$wsBinance.Range("B8").Value = "JTOUSDT"
$wsBinance.Range("C8").Value = 0.02210921954455465
$wsBinance.Range("E8").Value = 0.02211410880141988
$wsBinance.Range("B9").Value = "ETHUSDC"
$wsBinance.Range("C9").Value = 0.001330595263976798
$wsBinance.Range("D9").Value = -0.001330595263976798
$wsBinance.Range("B10").Value = "ETHUSDT"
$wsBinance.Range("C10").Value = 0.0004435750373607078
$wsBinance.Range("D10").Value = -0.0004435750373607078
$wsBinance.Range("E10").Value = 0
$wsBinance.Range("B11").Value = "RUNEUSDT"
$wsBinance.Range("C11").Value = 0.01681520094165687
$wsBinance.Range("E11").Value = 0.01681802892701537
$wsBinance.Range("B12").Value = "MATICUSDT"
$wsBinance.Range("C12").Value = 0.01161710037174593
$wsBinance.Range("D12").Value = 0
$wsBinance.Range("E12").Value = 0.01161845009875555
$wsBinance.Range("B13").Value = "DOTUSDT"
$wsBinance.Range("C13").Value = 0.01379310344828047
$wsBinance.Range("D13").Value = -0.04137931034482915
$wsBinance.Range("E13").Value = -0.02759001241550255
$wsBinance.Range("B14").Value = "ATOMUSDT"
$wsBinance.Range("C14").Value = 0.008942144326204469
$wsBinance.Range("D14").Value = -0.02682643297862929
$wsBinance.Range("E14").Value = -0.01788588803434688
$wsBinance.Range("B15").Value = "ETHBTC"
$wsBinance.Range("C15").Value = 0.01883239171375341
$wsBinance.Range("E15").Value = 0.0188359389715635
$wsBinance.Range("B16").Value = "LINKUSDT"
$wsBinance.Range("C16").Value = 0.00671546571754379
$wsBinance.Range("D16").Value = -0.00671546571754379
$wsBinance.Range("B17").Value = "DOGEUSDT"
$wsBinance.Range("C17").Value = 0.01049979000421042
$wsBinance.Range("D17").Value = -0.01049979000421042
$wsBinance.Range("B18").Value = "XRPUSDT"
$wsBinance.Range("C18").Value = 0.01601537475976762
$wsBinance.Range("D18").Value = -0.01601537475976762
$wsBinance.Range("B19").Value = "BTCUSDC"
$wsBinance.Range("C19").Value = 0.00004709817543269348
$wsBinance.Range("D19").Value = -0.007370864456390223
$wsBinance.Range("E19").Value = -0.007323769730319446
$wsBinance.Range("B20").Value = "AVAXUSDT"
$wsBinance.Range("C20").Value = 0.02495632642874472
$wsBinance.Range("D20").Value = 0
$wsBinance.Range("E20").Value = 0.0249625561657464
$wsBinance.Range("B21").Value = "INJUSDT"
$wsBinance.Range("C21").Value = 0.009631127805066338
$wsBinance.Range("D21").Value = -0.006420751870051829
$wsBinance.Range("E21").Value = 0.003210685160205707

Write-Output "done"
